$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.363.22'
$ws.Range('E2').Value = '  +0.03%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.801.96'
$ws.Range('E3').Value = '  +0.77%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.57'
$ws.Range('E5').Value = '  +0.76%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.575'
$ws.Range('E6').Value = '  +4.20%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.26'
$ws.Range('E8').Value = '  +11.20%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0694'
$ws.Range('E10').Value = '  +0.68%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +2.16%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.060.80'
$ws.Range('E12').Value = '  +0.67%  '

# Row 13
$ws.Range('E13').Value = '  +6.60%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.814.04'
$ws.Range('E14').Value = '  +1.44%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.644'
$ws.Range('E15').Value = '  +1.84%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.52'
$ws.Range('E16').Value = '  +6.03%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.307.34'
$ws.Range('E17').Value = '  -0.17%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.20'
$ws.Range('E18').Value = '  +1.42%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.66'
$ws.Range('E19').Value = '  +0.64%  '

# Row 20
$ws.Range('E20').Value = '  +0.05%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.73'
$ws.Range('E21').Value = '  +4.77%  '

# Row 22
$ws.Range('E22').Value = '  -0.03%  '

# Row 23
$ws.Range('E23').Value = '  +0.62%  '

# Row 24
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  +3.16%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.51'
$ws.Range('E25').Value = '  +3.40%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.95'
$ws.Range('E26').Value = '  +9.53%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.88'
$ws.Range('E27').Value = '  +2.55%  '

# Row 28
$ws.Range('E28').Value = '  +2.48%  '

# Row 29
$ws.Range('E29').Value = '  -0.08%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.01'
$ws.Range('E30').Value = '  +1.12%  '

# Row 32
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.83'
$ws.Range('E32').Value = '  +1.23%  '

# Row 33
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.24'
$ws.Range('E33').Value = '  +0.79%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').Value = '  +0.70%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.399.64'
$ws.Range('E35').Value = '  -0.07%  '

# Row 36
$ws.Range('E36').Value = '  -0.47%  '

# Row 37
$ws.Range('E37').Value = '  -4.10%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.06'
$ws.Range('E38').Value = '  +0.16%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0190'
$ws.Range('E39').Value = '  +0.24%  '

# Row 40
$ws.Range('E40').Value = '  +10.39%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.964'
$ws.Range('E41').Value = '  +3.35%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '82.53'
$ws.Range('E42').Value = '  -2.25%  '

# Row 43
$ws.Range('E43').Value = '  +0.46%  '

# Row 44
$ws.Range('E44').Value = '  +0.70%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.43'
$ws.Range('E45').Value = '  -2.72%  '

# Row 46
$ws.Range('E46').Value = '  -3.61%  '

# Row 47
$ws.Range('E47').Value = '  +0.81%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.963.02'
$ws.Range('E48').Value = '  +0.87%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.46'
$ws.Range('E49').Value = '  -0.14%  '

# Row 50
$ws.Range('E50').Value = '  +0.06%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0124'
$ws.Range('E51').Value = '  -3.18%  '
